# Set a shape's Left/Top/Width/Height precisely to a target EMU value,
# compensating for the single-precision (float32) storage PowerPoint uses
# for these properties (the engine truncates pts*12700 to EMU, so we search
# for a nearby point value whose float32 cast truncates to exactly the target).
function EmuToPt($emu) {
    $pts = $emu / 12700.0
    for ($i = 0; $i -le 5000; $i++) {
        foreach ($sign in @(1, -1)) {
            if ($i -eq 0 -and $sign -eq -1) { continue }
            $cand = $pts + ($sign * $i * 0.0000001)
            $f = [System.Single]$cand
            $back = [int64]([double]$f * 12700.0)
            if ($back -eq $emu) {
                return $cand
            }
        }
    }
    return $pts
}

function SetShapeEmu($sh, $offX, $offY, $extCx, $extCy) {
    if ($null -ne $offX) { $sh.Left = (EmuToPt $offX) }
    if ($null -ne $offY) { $sh.Top = (EmuToPt $offY) }
    if ($null -ne $extCx) { $sh.Width = (EmuToPt $extCx) }
    if ($null -ne $extCy) { $sh.Height = (EmuToPt $extCy) }
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

SetShapeEmu ($s.Shapes.Item("Rectangle 52")) $null $null 2474928 1787440
SetShapeEmu ($s.Shapes.Item("Rectangle 55")) $null $null 2474928 1797230
SetShapeEmu ($s.Shapes.Item("Rectangle 56")) $null $null 7741839 4075311
SetShapeEmu ($s.Shapes.Item("Rectangle 57")) 1467065 3639810 2470969 1520867
SetShapeEmu ($s.Shapes.Item("Rectangle 58")) $null $null 9585599 5157911
SetShapeEmu ($s.Shapes.Item("Rectangle 59")) 1300213 788789 2782099 4596012
SetShapeEmu ($s.Shapes.Item("Graphic 64")) 1468789 3648856 $null $null
SetShapeEmu ($s.Shapes.Item("Graphic 68")) 2444738 4141116 $null $null
SetShapeEmu ($s.Shapes.Item("TextBox 69")) 1923035 4637749 $null $null
SetShapeEmu ($s.Shapes.Item("Graphic 71")) 9180861 2812510 $null $null
SetShapeEmu ($s.Shapes.Item("TextBox 72")) 8640283 3394818 $null $null
SetShapeEmu ($s.Shapes.Item("TextBox 73")) 8640283 4640332 $null $null
SetShapeEmu ($s.Shapes.Item("Graphic 74")) 9180861 3993610 $null $null
SetShapeEmu ($s.Shapes.Item("TextBox 75")) 4136101 4648213 $null $null
SetShapeEmu ($s.Shapes.Item("Graphic 76")) 4539707 4057375 $null $null
SetShapeEmu ($s.Shapes.Item("Rectangle 79")) 5671319 3659964 2470969 1424609
SetShapeEmu ($s.Shapes.Item("Rectangle 80")) 5517334 798579 2782099 4586222
SetShapeEmu ($s.Shapes.Item("Graphic 82")) 5673043 3654589 $null $null
SetShapeEmu ($s.Shapes.Item("Graphic 85")) 6585447 4090316 $null $null
SetShapeEmu ($s.Shapes.Item("TextBox 86")) 6063744 4586949 $null $null
SetShapeEmu ($s.Shapes.Item("Straight Arrow Connector 89")) 5130741 4352892 $null $null
SetShapeEmu ($s.Shapes.Item("Straight Arrow Connector 90")) 3202660 4352892 $null $null
SetShapeEmu ($s.Shapes.Item("TextBox 4")) 9311797 120699 $null $null

# Remove the CIDR label textboxes (commit: "Removed CIDR numbers from architec diagram")
$s.Shapes.Item("TextBox 77").Delete()
$s.Shapes.Item("TextBox 78").Delete()
$s.Shapes.Item("TextBox 87").Delete()
$s.Shapes.Item("TextBox 88").Delete()
$s.Shapes.Item("TextBox 91").Delete()
